# Applies the updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.827.36'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '1.889.43'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('D4').Value = "'0.9998"
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'0.7776"
$ws.Range('E5').Value = '  -1.93%  '
$ws.Range('D6').Value = "'243.79"
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.95%  '
$ws.Range('D9').Value = "'0.07329"
$ws.Range('E9').Value = '  +4.24%  '
$ws.Range('D10').Value = "'25.28"
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').Value = "'0.08135"
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').Value = "'0.7648"
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = "'5.450"
$ws.Range('D14').Value = '1.880.07'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('D15').Value = "'92.98"
$ws.Range('E15').Value = '  +1.21%  '
$ws.Range('D16').Value = "'6.196"
$ws.Range('E16').Value = '  +4.52%  '
$ws.Range('D17').Value = '29.848.37'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('D19').Value = "'245.42"
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').Value = "'0.000007848"
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('D21').Value = "'0.9994"
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.146.99'
$ws.Range('E22').Value = '  +0.04%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = "'8.148"
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').Value = "'0.1577"
$ws.Range('E25').Value = '  -3.05%  '
$ws.Range('D26').Value = "'9.420"
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').Value = "'160.93"
$ws.Range('D28').Value = "'18.73"
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').Value = "'2.032"
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').Value = "'1.447"
$ws.Range('E30').Value = '  +5.74%  '
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('D32').Value = "'4.470"
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').Value = "'0.05582"
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('D34').Value = "'4.071"
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('E35').Value = '  -1.12%  '
$ws.Range('D36').Value = "'0.7536"
$ws.Range('E36').Value = '  +2.74%  '
$ws.Range('D37').Value = "'0.9969"
$ws.Range('E37').Value = '  -0.43%  '
$ws.Range('D38').Value = "'2.632"
$ws.Range('E38').Value = '  -3.06%  '
$ws.Range('D39').Value = "'0.01934"
$ws.Range('E39').Value = '  +0.83%  '
$ws.Range('D40').Value = "'2.777"
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('D41').Value = '1.138.36'
$ws.Range('E41').Value = '  +10.37%  '
$ws.Range('D42').Value = "'0.4439"
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('D43').Value = "'73.67"
$ws.Range('E43').Value = '  +2.45%  '
$ws.Range('D44').Value = "'5.956"
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').Value = "'0.8527"
$ws.Range('E45').Value = '  +1.51%  '
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').Value = "'1.896"
$ws.Range('E47').Value = '  +2.26%  '
$ws.Range('D48').Value = "'101.94"
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Value = "'3.102"
$ws.Range('E49').Value = '  +6.02%  '
$ws.Range('D50').Value = "'9.787"
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').Value = "'7.482"
$ws.Range('E51').Value = '  +0.68%  '
